$wb = $excel.ActiveWorkbook

# --- y_fitted_on_begin_2016 (sheet1) ---
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws.Cells.Item(2,1).Value2 = 1991
$ws.Cells.Item(2,2).Value2 = 16.34851261277244
$ws.Cells.Item(3,1).Value2 = 1992
$ws.Cells.Item(3,2).Value2 = 15.98158659889426
$ws.Cells.Item(4,1).Value2 = 1993
$ws.Cells.Item(4,2).Value2 = 15.78070533993826
$ws.Cells.Item(5,1).Value2 = 1994
$ws.Cells.Item(5,2).Value2 = 15.62618578580509
$ws.Cells.Item(6,1).Value2 = 1995
$ws.Cells.Item(6,2).Value2 = 16.13450302075714
$ws.Cells.Item(7,1).Value2 = 1996
$ws.Cells.Item(7,2).Value2 = 15.02320778762683
$ws.Cells.Item(8,1).Value2 = 1997
$ws.Cells.Item(8,2).Value2 = 14.60751480252285
$ws.Cells.Item(9,1).Value2 = 1998
$ws.Cells.Item(9,2).Value2 = 14.53307206947269
$ws.Cells.Item(10,1).Value2 = 1999
$ws.Cells.Item(10,2).Value2 = 13.96331715652476
$ws.Cells.Item(11,1).Value2 = 2000
$ws.Cells.Item(11,2).Value2 = 13.25523150715377
$ws.Cells.Item(12,1).Value2 = 2001
$ws.Cells.Item(12,2).Value2 = 12.76125226968025
$ws.Cells.Item(13,1).Value2 = 2002
$ws.Cells.Item(13,2).Value2 = 12.16293585454992
$ws.Cells.Item(14,1).Value2 = 2003
$ws.Cells.Item(14,2).Value2 = 11.81488314858483
$ws.Cells.Item(15,1).Value2 = 2004
$ws.Cells.Item(15,2).Value2 = 11.29679392207856
$ws.Cells.Item(16,1).Value2 = 2005
$ws.Cells.Item(16,2).Value2 = 10.79010221186542
$ws.Cells.Item(17,1).Value2 = 2006
$ws.Cells.Item(17,2).Value2 = 10.19286611920129
$ws.Cells.Item(18,1).Value2 = 2007
$ws.Cells.Item(18,2).Value2 = 9.658172001456801
$ws.Cells.Item(19,1).Value2 = 2008
$ws.Cells.Item(19,2).Value2 = 9.12388564631158
$ws.Cells.Item(20,1).Value2 = 2009
$ws.Cells.Item(20,2).Value2 = 9.127175872459807
$ws.Cells.Item(21,1).Value2 = 2010
$ws.Cells.Item(21,2).Value2 = 9.429053375060612
$ws.Cells.Item(22,1).Value2 = 2011
$ws.Cells.Item(22,2).Value2 = 9.31925680585835
$ws.Cells.Item(23,1).Value2 = 2012
$ws.Cells.Item(23,2).Value2 = 9.191184095254377
$ws.Cells.Item(24,1).Value2 = 2013
$ws.Cells.Item(24,2).Value2 = 9.286579300688176
$ws.Cells.Item(25,1).Value2 = 2014
$ws.Cells.Item(25,2).Value2 = 9.229008577443738
$ws.Cells.Item(26,1).Value2 = 2015
$ws.Cells.Item(26,2).Value2 = 9.277717568678062
$ws.Cells.Item(27,1).Value2 = 2016
$ws.Cells.Item(27,2).Value2 = 9.246561685550196

# --- y_pred_on_2017_2021 (sheet2) ---
$ws = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws.Cells.Item(2,1).Value2 = 2017
$ws.Cells.Item(2,2).Value2 = 9.091551032237028
$ws.Cells.Item(3,1).Value2 = 2018
$ws.Cells.Item(3,2).Value2 = 9.034868436721139
$ws.Cells.Item(4,1).Value2 = 2019
$ws.Cells.Item(4,2).Value2 = 8.872825948994585
$ws.Cells.Item(5,1).Value2 = 2020
$ws.Cells.Item(5,2).Value2 = 8.671859836956372
$ws.Cells.Item(6,1).Value2 = 2021
$ws.Cells.Item(6,2).Value2 = 8.455512845908467

# --- y_fitted_on_begin_2021 (sheet3) ---
$ws = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws.Rows.Item(33).Delete()
$ws.Cells.Item(2,1).Value2 = 1991
$ws.Cells.Item(2,2).Value2 = 16.35780381523324
$ws.Cells.Item(3,1).Value2 = 1992
$ws.Cells.Item(3,2).Value2 = 15.99708723071187
$ws.Cells.Item(4,1).Value2 = 1993
$ws.Cells.Item(4,2).Value2 = 15.80822708564829
$ws.Cells.Item(5,1).Value2 = 1994
$ws.Cells.Item(5,2).Value2 = 15.61419812433326
$ws.Cells.Item(6,1).Value2 = 1995
$ws.Cells.Item(6,2).Value2 = 16.15905137132411
$ws.Cells.Item(7,1).Value2 = 1996
$ws.Cells.Item(7,2).Value2 = 14.98637656597831
$ws.Cells.Item(8,1).Value2 = 1997
$ws.Cells.Item(8,2).Value2 = 14.60131971262905
$ws.Cells.Item(9,1).Value2 = 1998
$ws.Cells.Item(9,2).Value2 = 14.50597442511093
$ws.Cells.Item(10,1).Value2 = 1999
$ws.Cells.Item(10,2).Value2 = 13.94005769589224
$ws.Cells.Item(11,1).Value2 = 2000
$ws.Cells.Item(11,2).Value2 = 13.29492189885921
$ws.Cells.Item(12,1).Value2 = 2001
$ws.Cells.Item(12,2).Value2 = 12.77961544559738
$ws.Cells.Item(13,1).Value2 = 2002
$ws.Cells.Item(13,2).Value2 = 12.13474002445686
$ws.Cells.Item(14,1).Value2 = 2003
$ws.Cells.Item(14,2).Value2 = 11.81139876245226
$ws.Cells.Item(15,1).Value2 = 2004
$ws.Cells.Item(15,2).Value2 = 11.27988034478522
$ws.Cells.Item(16,1).Value2 = 2005
$ws.Cells.Item(16,2).Value2 = 10.7605353213657
$ws.Cells.Item(17,1).Value2 = 2006
$ws.Cells.Item(17,2).Value2 = 10.18235583274619
$ws.Cells.Item(18,1).Value2 = 2007
$ws.Cells.Item(18,2).Value2 = 9.678251857101811
$ws.Cells.Item(19,1).Value2 = 2008
$ws.Cells.Item(19,2).Value2 = 9.177854239585541
$ws.Cells.Item(20,1).Value2 = 2009
$ws.Cells.Item(20,2).Value2 = 9.163506200746166
$ws.Cells.Item(21,1).Value2 = 2010
$ws.Cells.Item(21,2).Value2 = 9.398282865931689
$ws.Cells.Item(22,1).Value2 = 2011
$ws.Cells.Item(22,2).Value2 = 9.332768226404772
$ws.Cells.Item(23,1).Value2 = 2012
$ws.Cells.Item(23,2).Value2 = 9.200742868701889
$ws.Cells.Item(24,1).Value2 = 2013
$ws.Cells.Item(24,2).Value2 = 9.236044495030399
$ws.Cells.Item(25,1).Value2 = 2014
$ws.Cells.Item(25,2).Value2 = 9.228445665686619
$ws.Cells.Item(26,1).Value2 = 2015
$ws.Cells.Item(26,2).Value2 = 9.266104444289955
$ws.Cells.Item(27,1).Value2 = 2016
$ws.Cells.Item(27,2).Value2 = 9.192078825911381
$ws.Cells.Item(28,1).Value2 = 2017
$ws.Cells.Item(28,2).Value2 = 8.99858008755346
$ws.Cells.Item(29,1).Value2 = 2018
$ws.Cells.Item(29,2).Value2 = 8.925955819152035
$ws.Cells.Item(30,1).Value2 = 2019
$ws.Cells.Item(30,2).Value2 = 8.777554871029992
$ws.Cells.Item(31,1).Value2 = 2020
$ws.Cells.Item(31,2).Value2 = 8.6486165104794
$ws.Cells.Item(32,1).Value2 = 2021
$ws.Cells.Item(32,2).Value2 = 8.861079728764233

# --- y_pred_on_2022_2026 (sheet4) ---
$ws = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws.Cells.Item(2,1).Value2 = 2022
$ws.Cells.Item(2,2).Value2 = 8.643679258073417
$ws.Cells.Item(3,1).Value2 = 2023
$ws.Cells.Item(3,2).Value2 = 8.754281420825849
$ws.Cells.Item(4,1).Value2 = 2024
$ws.Cells.Item(4,2).Value2 = 8.817577466088505
$ws.Cells.Item(5,1).Value2 = 2025
$ws.Cells.Item(5,2).Value2 = 8.908767729247709
$ws.Cells.Item(6,1).Value2 = 2026
$ws.Cells.Item(6,2).Value2 = 8.994035717411041
